$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.913.34'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '1.816.04'
$ws.Range("E3").Value = '  +0.80%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.07'
$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4649'
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3665'
$ws.Range("E8").Value = '  -0.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07362'
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8705'
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.31'
$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("D12").Value = '1.851.34'
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("E13").Value = '  +0.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07088'
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.508'
$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.52'
$ws.Range("E16").Value = '  -0.86%  '

$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008718'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("D21").Value = '26.944.00'
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.302'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").Value = '2.070.25'
$ws.Range("E24").Value = '  +0.39%  '

$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("E26").Value = '  -0.50%  '

$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.132'
$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.271'
$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.32'
$ws.Range("E30").Value = '  -0.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08895'
$ws.Range("E31").Value = '  -0.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7562'
$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("E33").Value = '  +0.78%  '

$ws.Range("E34").Value = '  +0.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.907'
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.084'
$ws.Range("E37").Value = '  -1.32%  '

$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01948'
$ws.Range("E39").Value = '  -0.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.985'
$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.237'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5296'
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.306'
$ws.Range("E43").Value = '  -3.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1654'
$ws.Range("E44").Value = '  -0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.442'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4866'
$ws.Range("E46").Value = '  -2.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.45'
$ws.Range("E47").Value = '  +1.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.38'
$ws.Range("E49").Value = '  -0.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.661'
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("E51").Value = '  +0.07%  '

